# Add a new log entry (row 4) to the lab-notebook sheet:
#   - A4: timestamp of the new observation (same date/time style as A3)
#   - E4: the new "Observations" note, word-wrapped like the rest of column E
# Also retro-fits the header cell E1 ("Observations") with the wrap-text
# style already used by the rest of the table, and moves the selection to
# the freshly-entered cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell E1: give it the same wrap-text formatting as D1/etc. ---
$ws.Range("E1").WrapText = $true

# --- New row 4 ---
# A4 needs the exact same date/time format as the other timestamp cells
# (A2/A3). Copy A3's formatting first, then overwrite just the value so we
# reuse the existing style instead of minting a new numFmt.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 43401.816666666666
$ws.Application.CutCopyMode = $false

# E4 holds the new observation text; wrap it like the other Observations/
# Description cells.
$ws.Range("E4").Value = "Not able to get CodeBaseRegister to manage memory without memory exceptions.  Will abandon this work for now, undo the CodeBaseRegister, and continue this branch to research other options."
$ws.Range("E4").WrapText = $true

# Match the row height the wrapped text occupies in the source workbook.
$ws.Rows(4).RowHeight = 28.8

# Leave the selection on the newly added cell, as in the saved workbook.
$ws.Range("E4").Select()
